$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 358, shifting existing rows 358:446 down to 359:447
$ws.Rows.Item(358).Insert()

# Populate the newly inserted row 358 with the new record's data
$ws.Cells.Item(358, 1).Value = 5
$ws.Cells.Item(358, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(358, 3).Value = "Maule"
$ws.Cells.Item(358, 4).Value = 44551
$ws.Cells.Item(358, 5).Value = 7
$ws.Cells.Item(358, 6).Value = 100112020
$ws.Cells.Item(358, 7).Value = "Tomate"
$ws.Cells.Item(358, 8).Value = "Larga vida"
$ws.Cells.Item(358, 9).Value = "Primera"
$ws.Cells.Item(358, 10).Value = 4000
$ws.Cells.Item(358, 11).Value = 10000
$ws.Cells.Item(358, 12).Value = 10000
$ws.Cells.Item(358, 13).Value = 10000
$ws.Cells.Item(358, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(358, 15).Value = "Región del Maule"
$ws.Cells.Item(358, 16).Value = 667
$ws.Cells.Item(358, 17).Value = 15
$ws.Cells.Item(358, 18).Value = "Hortaliza"
